$wb = $excel.ActiveWorkbook

$wsMacro = $wb.Worksheets.Item("Macro_taxonomy")
$wsCosts = $wb.Worksheets.Item("Costs")

# --- Insert a new row at row 16 (Other/Urban ME+MEO/LWAL 0.5), pushing the
# existing Rural block (old rows 16-24) down to rows 17-25 ---------------
$wsMacro.Rows.Item(16).Insert()

# Row 15 (Other / Urban / MATO) proportion changes from 1 -> 0.5
$wsMacro.Range("D15").Value = 0.5

# New row 16: Other / Urban / ME+MEO/LWAL / 0.5
$wsMacro.Range("A16").Value = "Other"
$wsMacro.Range("B16").Value = "Urban"
$wsMacro.Range("C16").Value = "ME+MEO/LWAL"
$wsMacro.Range("D16").Value = 0.5

# Old row 24 (now row 25, Other / Rural / MATO) proportion changes from 1 -> 0.5
$wsMacro.Range("D25").Value = 0.5

# --- Insert two new rows after the shifted Rural block (now ending at row 25)
# to hold the two new Other/Rural entries -> rows 26 and 27 --------------
$wsMacro.Rows.Item(26).Insert()
$wsMacro.Rows.Item(27).Insert()

$wsMacro.Range("A26").Value = "Other"
$wsMacro.Range("B26").Value = "Rural"
$wsMacro.Range("C26").Value = "EWV/LN"
$wsMacro.Range("D26").Value = 0.25

$wsMacro.Range("A27").Value = "Other"
$wsMacro.Range("B27").Value = "Rural"
$wsMacro.Range("C27").Value = "ME+MEO/LWAL"
$wsMacro.Range("D27").Value = 0.25

# --- Insert a blank placeholder row at (original) row 44, pushing the rest
# of the blank placeholder rows down by one and adding a new trailing blank
# row at the bottom of the sheet. Three rows were already inserted above
# (at rows 16, 26 and 27), so original row 44 now sits at row 47. ---------
$wsMacro.Rows.Item(47).Insert()

# Active cell / selection on the Macro_taxonomy sheet
$wsMacro.Activate()
$wsMacro.Range("D17").Select()

# Costs sheet is no longer the active tab/sheet; its own stored selection
# stays as-is (still E1:E33) since only sheet activation changed.
$wsCosts.Range("E1:E33").Select()
$wsMacro.Activate()
$wsMacro.Range("D17").Select()
